$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# F2 currently holds the text "False" as a shared string.
# Push an actual boolean value (FALSE) into the cell instead.
$ws.Range("F2").Value = $false
